$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 71,3
$data[0,0] = 34.12942556846704
$data[0,1] = 29.31856681616817
$data[0,2] = 39.19264274160815
$data[1,0] = 26.1085361391383
$data[1,1] = 9.346835616438318
$data[1,2] = 43.30291095890406
$data[2,0] = 25.61246420446107
$data[2,1] = 8.127516438356157
$data[2,2] = 44.10502054794514
$data[3,0] = 32.05133063292968
$data[3,1] = 8.661695890410936
$data[3,2] = 52.71659452054792
$data[4,0] = 37.42199821876793
$data[4,1] = 24.32837267390138
$data[4,2] = 49.77000765492491
$data[5,0] = 29.02369712323709
$data[5,1] = 24.20142678428143
$data[5,2] = 34.00034884690801
$data[6,0] = 19.74371943946466
$data[6,1] = 16.68357572480302
$data[6,2] = 22.69771929828047
$data[7,0] = 34.17500557534638
$data[7,1] = 21.46197427701669
$data[7,2] = 46.87899118306347
$data[8,0] = 21.86237345406464
$data[8,1] = 18.52910375391417
$data[8,2] = 25.24287847222094
$data[9,0] = 31.0936790550825
$data[9,1] = 26.65696929328635
$data[9,2] = 35.61371071835922
$data[10,0] = 41.48439171289267
$data[10,1] = 34.68616755797021
$data[10,2] = 48.77938802541948
$data[11,0] = 28.33678478820592
$data[11,1] = 24.59742474637492
$data[11,2] = 32.39165826030511
$data[12,0] = 28.46121568602972
$data[12,1] = 24.06938743787805
$data[12,2] = 32.32682432831723
$data[13,0] = 39.01139191436205
$data[13,1] = 11.99314520547943
$data[13,2] = 62.50157385844747
$data[14,0] = 37.62268936181368
$data[14,1] = 24.50282853387094
$data[14,2] = 49.91878752935411
$data[15,0] = 41.614971646052
$data[15,1] = 34.73326651471679
$data[15,2] = 48.84536074197544
$data[16,0] = 25.90823312998291
$data[16,1] = 7.715773972602733
$data[16,2] = 43.29995890410951
$data[17,0] = 35.37725806911243
$data[17,1] = 31.06512943402012
$data[17,2] = 40.0252616163414
$data[18,0] = 42.60071102594345
$data[18,1] = 29.73003814959765
$data[18,2] = 53.29272130324579
$data[19,0] = 38.97228270249702
$data[19,1] = 18.35535730593601
$data[19,2] = 57.66375536529679
$data[20,0] = 23.8225290995534
$data[20,1] = 20.27041489334389
$data[20,2] = 27.42934697190262
$data[21,0] = 38.23173071104372
$data[21,1] = 9.725656073059337
$data[21,2] = 63.4669600913242
$data[22,0] = 37.55191223670519
$data[22,1] = 24.4524587021684
$data[22,2] = 49.90610727821257
$data[23,0] = 30.869641371766
$data[23,1] = 18.67872104370511
$data[23,2] = 43.02911958904106
$data[24,0] = 38.51791641056233
$data[24,1] = 17.09056844748853
$data[24,2] = 57.20145093607307
$data[25,0] = 31.91699371125829
$data[25,1] = 27.64797658862376
$data[25,2] = 36.19512596740906
$data[26,0] = 31.57688499601382
$data[26,1] = 27.01270754546501
$data[26,2] = 35.98355368103985
$data[27,0] = 42.07436530901038
$data[27,1] = 14.51732054794519
$data[27,2] = 67.95699999999989
$data[28,0] = 38.56013081978724
$data[28,1] = 31.59071560864117
$data[28,2] = 44.92904038547671
$data[29,0] = 26.3850687612412
$data[29,1] = 22.63395273274789
$data[29,2] = 30.07027319180291
$data[30,0] = 22.73950625421337
$data[30,1] = 19.00281697755642
$data[30,2] = 26.21903017727035
$data[31,0] = 28.53547808686232
$data[31,1] = 24.45439953724479
$data[31,2] = 32.4557622560317
$data[32,0] = 41.13276909015847
$data[32,1] = 35.13540420259539
$data[32,2] = 46.8216952986857
$data[33,0] = 39.88184167888606
$data[33,1] = 27.27994400403248
$data[33,2] = 51.5019811308564
$data[34,0] = 29.1403884615348
$data[34,1] = 24.53153260955471
$data[34,2] = 33.70303633282042
$data[35,0] = 41.0880898100937
$data[35,1] = 19.20888748858443
$data[35,2] = 59.68280545662101
$data[36,0] = 31.8772654546153
$data[36,1] = 27.55649750432548
$data[36,2] = 36.61292671085383
$data[37,0] = 44.0651516130504
$data[37,1] = 36.94323515842868
$data[37,2] = 51.38939509448808
$data[38,0] = 32.34747513956009
$data[38,1] = 12.10359269406385
$data[38,2] = 50.57102511415522
$data[39,0] = 26.32256772747394
$data[39,1] = 21.45186075426071
$data[39,2] = 31.26842164047619
$data[40,0] = 33.04668048007434
$data[40,1] = 28.59460170512763
$data[40,2] = 37.69161692963961
$data[41,0] = 38.55414186684006
$data[41,1] = 31.60815511478123
$data[41,2] = 44.90649165056023
$data[42,0] = 24.26741289046595
$data[42,1] = 20.60921519300247
$data[42,2] = 28.26519458625814
$data[43,0] = 36.96422605449266
$data[43,1] = 31.95345325939918
$data[43,2] = 42.30085496088262
$data[44,0] = 24.24853275888623
$data[44,1] = 20.79147969320045
$data[44,2] = 27.78348250439363
$data[45,0] = 25.84286615829155
$data[45,1] = 21.34482895753736
$data[45,2] = 29.9041903708369
$data[46,0] = 41.6130932237929
$data[46,1] = 34.82354814430418
$data[46,2] = 48.90366767460502
$data[47,0] = 30.1781930606387
$data[47,1] = 8.454687671232865
$data[47,2] = 47.00985068493145
$data[48,0] = 41.2550864366804
$data[48,1] = 19.77558630136981
$data[48,2] = 59.72596283105022
$data[49,0] = 25.89429457946789
$data[49,1] = 7.294867123287662
$data[49,2] = 43.46526849315065
$data[50,0] = 22.10038745071709
$data[50,1] = 18.84506979246748
$data[50,2] = 25.4421951089014
$data[51,0] = 44.13256122762554
$data[51,1] = 37.15316099798428
$data[51,2] = 51.37792930276139
$data[52,0] = 38.27178301174649
$data[52,1] = 10.9088625570776
$data[52,2] = 62.19631506849313
$data[53,0] = 39.03601558699202
$data[53,1] = 16.83902945205475
$data[53,2] = 56.99191906392691
$data[54,0] = 26.3285266134616
$data[54,1] = 22.58291848815886
$data[54,2] = 30.29180453623434
$data[55,0] = 22.74383543772412
$data[55,1] = 19.01217846172724
$data[55,2] = 26.23681375755668
$data[56,0] = 40.59902330730518
$data[56,1] = 17.63946328767119
$data[56,2] = 58.9572366210046
$data[57,0] = 28.77200455725475
$data[57,1] = 24.36009829744212
$data[57,2] = 33.22746851960453
$data[58,0] = 39.4692762305454
$data[58,1] = 11.72526712328767
$data[58,2] = 64.62098207762558
$data[59,0] = 30.02478003034906
$data[59,1] = 25.90193449516411
$data[59,2] = 34.40939345051892
$data[60,0] = 31.02341694954274
$data[60,1] = 26.59638861372367
$data[60,2] = 35.54118692945911
$data[61,0] = 34.05092947585607
$data[61,1] = 29.84357358760177
$data[61,2] = 37.93909210361434
$data[62,0] = 39.22473083552515
$data[62,1] = 11.08142420091323
$data[62,2] = 64.50699486301372
$data[63,0] = 20.69713323799546
$data[63,1] = 17.42888855813554
$data[63,2] = 24.02051526496971
$data[64,0] = 38.81281171334494
$data[64,1] = 33.72934768000191
$data[64,2] = 44.19705028459467
$data[65,0] = 36.23189411997058
$data[65,1] = 31.21503068953969
$data[65,2] = 41.41321923989365
$data[66,0] = 35.69908291285977
$data[66,1] = 16.42873607305933
$data[66,2] = 53.80674406392691
$data[67,0] = 25.78019013335827
$data[67,1] = 22.37234408637107
$data[67,2] = 29.55053044746207
$data[68,0] = 32.3372421531549
$data[68,1] = 12.80817899543372
$data[68,2] = 49.99497305936065
$data[69,0] = 27.45792606516065
$data[69,1] = 18.84543450823921
$data[69,2] = 34.6010936641518
$data[70,0] = 33.45446971398252
$data[70,1] = 26.02234158500374
$data[70,2] = 40.56987031061694

$ws.Range("A2:C72").Value = $data
